$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Update the division-problem cells in the practice table.
# Each cell's old expression is replaced with its new expression in place,
# using wdReplaceOne (1) scoped to that single cell's Range so that
# duplicate expressions elsewhere in the table are not touched.

# Row 1
$t.Cell(1,1).Range.Find.Execute("42÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷5=", 1)
$t.Cell(1,2).Range.Find.Execute("70÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷3=", 1)
$t.Cell(1,3).Range.Find.Execute("42÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷8=", 1)
$t.Cell(1,4).Range.Find.Execute("62÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "91÷7=", 1)
$t.Cell(1,5).Range.Find.Execute("65÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷5=", 1)

# Row 5
$t.Cell(5,1).Range.Find.Execute("44÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷8=", 1)
$t.Cell(5,2).Range.Find.Execute("58÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷8=", 1)
$t.Cell(5,3).Range.Find.Execute("89÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷5=", 1)
$t.Cell(5,4).Range.Find.Execute("76÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷3=", 1)
$t.Cell(5,5).Range.Find.Execute("18÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷3=", 1)

# Row 9
$t.Cell(9,1).Range.Find.Execute("82÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷7=", 1)
$t.Cell(9,2).Range.Find.Execute("86÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷5=", 1)
$t.Cell(9,3).Range.Find.Execute("14÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷7=", 1)
$t.Cell(9,4).Range.Find.Execute("89÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷9=", 1)
$t.Cell(9,5).Range.Find.Execute("81÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷7=", 1)

# Row 13
$t.Cell(13,1).Range.Find.Execute("66÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷6=", 1)
$t.Cell(13,2).Range.Find.Execute("60÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷6=", 1)
$t.Cell(13,3).Range.Find.Execute("44÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "65÷4=", 1)
$t.Cell(13,4).Range.Find.Execute("46÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷6=", 1)
$t.Cell(13,5).Range.Find.Execute("83÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷7=", 1)

# Row 17
$t.Cell(17,1).Range.Find.Execute("61÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷5=", 1)
$t.Cell(17,2).Range.Find.Execute("58÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷8=", 1)
$t.Cell(17,3).Range.Find.Execute("90÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "89÷9=", 1)
$t.Cell(17,4).Range.Find.Execute("22÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷6=", 1)
$t.Cell(17,5).Range.Find.Execute("28÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "94÷2=", 1)
